$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.079.59'
$ws.Range('E2').Value = '  -4.43%  '

# Row 3
$ws.Range('D3').Value = '3.692.55'
$ws.Range('E3').Value = '  -5.02%  '

# Row 4
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.11'
$ws.Range('E5').Value = '  -0.46%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.74'
$ws.Range('E6').Value = '  +7.08%  '

# Row 7
$ws.Range('D7').Value = '3.680.12'
$ws.Range('E7').Value = '  -5.23%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.626'
$ws.Range('E8').Value = '  -7.19%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.08%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.713'
$ws.Range('E10').Value = '  -5.69%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.162'
$ws.Range('E11').Value = '  -8.37%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '55.99'
$ws.Range('E12').Value = '  +3.76%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000289'
$ws.Range('E13').Value = '  -10.77%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.35'
$ws.Range('E14').Value = '  -9.77%  '

# Row 15
$ws.Range('D15').Value = '4.290.01'
$ws.Range('E15').Value = '  -4.91%  '

# Row 16
$ws.Range('D16').Value = '3.703.14'
$ws.Range('E16').Value = '  -4.85%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.36'
$ws.Range('E17').Value = '  -7.83%  '

# Row 18
$ws.Range('E18').Value = '  -2.38%  '

# Row 19
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.12'
$ws.Range('E19').Value = '  -7.68%  '

# Row 20
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.78'
$ws.Range('E20').Value = '  -8.54%  '

# Row 21
$ws.Range('D21').Value = '67.937.16'
$ws.Range('E21').Value = '  -4.49%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '408.60'
$ws.Range('E22').Value = '  -6.84%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.55'
$ws.Range('E23').Value = '  -4.38%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.41'
$ws.Range('E24').Value = '  -6.95%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.01'
$ws.Range('E25').Value = '  -8.89%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.75'
$ws.Range('E26').Value = '  -8.37%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.01'
$ws.Range('E27').Value = '  -2.79%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.85'
$ws.Range('E28').Value = '  -6.41%  '

# Row 29
$ws.Range('E29').Value = '  +2.19%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.45'
$ws.Range('E30').Value = '  -9.17%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.66'
$ws.Range('E31').Value = '  -7.45%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.24'
$ws.Range('E32').Value = '  -11.00%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.46'
$ws.Range('E33').Value = '  -9.00%  '

# Row 34
$ws.Range('E34').Value = '  -7.27%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '43.39'
$ws.Range('E35').Value = '  -12.57%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.05'
$ws.Range('E36').Value = '  -8.95%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '590.61'
$ws.Range('E37').Value = '  -7.24%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0880'
$ws.Range('E38').Value = '  -11.54%  '

# Row 39
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.00%  '

# Row 40
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.399'
$ws.Range('E40').Value = '  -6.43%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.07%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.136'
$ws.Range('E42').Value = '  -6.22%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.75'
$ws.Range('E43').Value = '  -0.07%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.99'
$ws.Range('E44').Value = '  -8.93%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0435'
$ws.Range('E45').Value = '  -8.28%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.87'
$ws.Range('E46').Value = '  -14.85%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.15'
$ws.Range('E47').Value = '  -10.13%  '

# Row 48
$ws.Range('E48').Value = '  -4.02%  '

# Row 49
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.762.80'
$ws.Range('E49').Value = '  -2.60%  '

# Row 50
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.134'
$ws.Range('E50').Value = '  -7.42%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.13'
$ws.Range('E51').Value = '  -4.98%  '
